$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.347.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.877.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7200'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08017'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3136'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08155'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.871.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '94.70'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.221'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7099'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('E16').Value = '  +4.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008469'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.343.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.124.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.735'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1605'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.038'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.504'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.404'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.281'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.215'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05349'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.933'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7606'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.177'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01867'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.267.65'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.762'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.439'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '113.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9064'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '74.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000130'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.022.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.799'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5195'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.479'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4338'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
